$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$newRow = $t.Rows.Add()
$idx = $newRow.Index

$t.Cell($idx, 1).Range.Text = "7"
$t.Cell($idx, 2).Range.Text = "-"
$t.Cell($idx, 3).Range.Text = "0.377"
